$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (C, D, H) ---
# Stored column width = ColumnWidth + 5/6, so subtract that offset to land on the
# exact target widths of 60 / 50 / 33 characters.
$ws.Columns.Item(3).ColumnWidth = 60 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 50 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 33 - 0.8333333333333334

# Keep opportunity IDs (column A) as text, not auto-converted to numbers.
$ws.Range("A2:A14").NumberFormat = "@"

# Row 2: PwC Global Partnership
$ws.Range("A2").Value = "1328588"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328588"
$ws.Range("C2").Value = "My Way​ Operations & Innovation Coordinator​ 2026-2027"
$ws.Range("D2").Value = "40 Düsseldorf, Germany"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "146 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "PwC Global Partnership"
$ws.Range("E2").Interior.Color = 65535

# Row 3: Neospace A. I. Technologies
$ws.Range("A3").Value = "1331173"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1331173"
$ws.Range("C3").Value = "[Impact Brazil] - A. I. Technologies Developer Intern"
$ws.Range("D3").Value = "Uberlândia, MG, Brasil"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "1 applicant"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Neospace A. I. Technologies"

# Row 4: OK Lucky Star
$ws.Range("A4").Value = "1331172"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1331172"
$ws.Range("C4").Value = "Marketing Communication Intern/Coordinator"
$ws.Range("D4").Value = "Belgrade, Serbia"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "OK Lucky Star"

# Row 5: Azayeiz Football Academy
$ws.Range("A5").Value = "1331137"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1331137"
$ws.Range("C5").Value = "Sales Agent"
$ws.Range("D5").Value = "Sousse, Tunisie"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "Azayeiz Football Academy"

# Row 6: Lovely Professional University
$ws.Range("A6").Value = "1331115"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1331115"
$ws.Range("C6").Value = "Artificial Intelligence Intern"
$ws.Range("D6").Value = "Jalandhar, Punjab, India"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "Lovely Professional University"

# Row 7: Business 360
$ws.Range("A7").Value = "1331093"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1331093"
$ws.Range("C7").Value = "Digital Marketer"
$ws.Range("D7").Value = "Sousse, Tunisie"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "0 applicants"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "Business 360"

# Row 8: DENSsolutions B.V.
$ws.Range("A8").Value = "1330746"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1330746"
$ws.Range("C8").Value = "Software engineer (EU only)"
$ws.Range("D8").Value = "Delft, Nederland"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "104 applicants"
$ws.Range("G8").Value = "6 - 18 Months"
$ws.Range("H8").Value = "DENSsolutions B.V."

# Row 9: GSK - GlaxoSmithKline
$ws.Range("A9").Value = "1330676"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1330676"
$ws.Range("C9").Value = "General Medicines Marketing Analyst"
$ws.Range("D9").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "23 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "GSK - GlaxoSmithKline"

# Row 10: Nidec Global Appliance
$ws.Range("A10").Value = "1329526"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1329526"
$ws.Range("C10").Value = "[Impact Brazil] - Supplier Development Procurement Intern"
$ws.Range("D10").Value = "Joinville - Pirabeiraba, Joinville - SC, Brasil"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "36 applicants"
$ws.Range("G10").Value = "3 - 6 Months"
$ws.Range("H10").Value = "Nidec Global Appliance"

# Row 11: Nidec Global Appliance
$ws.Range("A11").Value = "1329430"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1329430"
$ws.Range("C11").Value = "[Impact Brazil] - Research and Development Lab Intern"
$ws.Range("D11").Value = "Joinville - Pirabeiraba, Joinville - SC, Brasil"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "58 applicants"
$ws.Range("G11").Value = "3 - 6 Months"
$ws.Range("H11").Value = "Nidec Global Appliance"

# Row 12: Tepma MEP Design
$ws.Range("A12").Value = "1328980"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1328980"
$ws.Range("C12").Value = "Engineering Intern"
$ws.Range("D12").Value = "Belgrade, Serbia"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "70 applicants"
$ws.Range("G12").Value = "9 - 12 Weeks"
$ws.Range("H12").Value = "Tepma MEP Design"

# Row 13: OK Roda
$ws.Range("A13").Value = "1328974"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1328974"
$ws.Range("C13").Value = "Youth Volleyball Assistant Coach"
$ws.Range("D13").Value = "Belgrade, Serbia"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "33 applicants"
$ws.Range("G13").Value = "9 - 12 Weeks"
$ws.Range("H13").Value = "OK Roda"

# Row 14: Ikan Experience
$ws.Range("A14").Value = "1316788"
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1316788"
$ws.Range("C14").Value = "Travel Coordinator"
$ws.Range("D14").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "135 applicants"
$ws.Range("G14").Value = "6 - 18 Months"
$ws.Range("H14").Value = "Ikan Experience"
